$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF); copy H1's formatting (bold, border,
# centered/top alignment) onto them first, then set the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I2:J44
$iValues = @(9,8,6,8,7,6,5,6,5,7,7,6,8,8,6,8,7,1,7,8,7,8,6,9,8,8,6,7,4,8,7,6,7,8,9,7,6,6,7,8,1,1,8)
$jValues = @(9,8,6,8,7,6,6,6,6,7,7,6,8,8,7,8,7,1,8,8,8,8,6,9,8,8,6,7,5,8,7,6,8,8,9,7,6,7,9,8,4,3,8)

for ($r = 2; $r -le 44; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
